$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet tab
$ws.Name = "repayment_20250901_20250916 (1)"

$ws.Cells.Item(2,4).Value = 45  # D2: 46 -> 45
$ws.Cells.Item(2,5).Value = "34,748,778.00"  # E2: '34,486,339.00' -> '34,748,778.00'
$ws.Cells.Item(2,6).Value = "341,892,945.00"  # F2: '332,666,040.00' -> '341,892,945.00'
$ws.Cells.Item(2,7).Value = "10.16"  # G2: '10.37' -> '10.16'
$ws.Cells.Item(2,8).Value = 15.519  # H2: 15.304 -> 15.519
$ws.Cells.Item(2,9).Value = 266  # I2: 252 -> 266
$ws.Cells.Item(2,11).Value = "9.72"  # K2: '10.13' -> '9.72'
$ws.Cells.Item(2,12).Value = "6.02"  # L2: '6.35' -> '6.02'

$ws.Cells.Item(3,4).Value = 41  # D3: 42 -> 41
$ws.Cells.Item(3,6).Value = "335,755,173.00"  # F3: '325,976,616.00' -> '335,755,173.00'
$ws.Cells.Item(3,7).Value = "9.24"  # G3: '9.52' -> '9.24'
$ws.Cells.Item(3,8).Value = 20.546  # H3: 20.086 -> 20.546
$ws.Cells.Item(3,9).Value = 266  # I3: 252 -> 266
$ws.Cells.Item(3,11).Value = "2.43"  # K3: '2.54' -> '2.43'
$ws.Cells.Item(3,12).Value = "3.01"  # L3: '3.17' -> '3.01'

$ws.Cells.Item(4,4).Value = 39  # D4: 36 -> 39
$ws.Cells.Item(4,5).Value = "23,648,532.00"  # E4: '22,568,848.00' -> '23,648,532.00'
$ws.Cells.Item(4,6).Value = "354,114,049.00"  # F4: '344,627,951.00' -> '354,114,049.00'
$ws.Cells.Item(4,7).Value = "6.68"  # G4: '6.55' -> '6.68'
$ws.Cells.Item(4,8).Value = 12.352  # H4: 11.498 -> 12.352
$ws.Cells.Item(4,9).Value = 268  # I4: 253 -> 268
$ws.Cells.Item(4,10).Value = 17  # J4: 16 -> 17
$ws.Cells.Item(4,11).Value = "5.30"  # K4: '5.37' -> '5.30'
$ws.Cells.Item(4,12).Value = "6.34"  # L4: '6.32' -> '6.34'

$ws.Cells.Item(5,4).Value = 41  # D5: 46 -> 41
$ws.Cells.Item(5,5).Value = "27,597,133.00"  # E5: '26,609,103.00' -> '27,597,133.00'
$ws.Cells.Item(5,6).Value = "315,737,797.00"  # F5: '307,121,314.00' -> '315,737,797.00'
$ws.Cells.Item(5,7).Value = "8.74"  # G5: '8.66' -> '8.74'
$ws.Cells.Item(5,8).Value = 20.452  # H5: 20.208 -> 20.452
$ws.Cells.Item(5,9).Value = 236  # I5: 222 -> 236
$ws.Cells.Item(5,10).Value = 12  # J5: 10 -> 12
$ws.Cells.Item(5,11).Value = "5.23"  # K5: '3.95' -> '5.23'
$ws.Cells.Item(5,12).Value = "5.08"  # L5: '4.50' -> '5.08'

$ws.Cells.Item(6,4).Value = 39  # D6: 40 -> 39
$ws.Cells.Item(6,5).Value = "26,306,910.00"  # E6: '25,792,031.00' -> '26,306,910.00'
$ws.Cells.Item(6,6).Value = "324,518,843.00"  # F6: '316,238,350.00' -> '324,518,843.00'
$ws.Cells.Item(6,7).Value = "8.11"  # G6: '8.16' -> '8.11'
$ws.Cells.Item(6,8).Value = 14.304  # H6: 14.285 -> 14.304
$ws.Cells.Item(6,9).Value = 266  # I6: 252 -> 266
$ws.Cells.Item(6,10).Value = 8  # J6: 7 -> 8
$ws.Cells.Item(6,11).Value = "2.80"  # K6: '2.66' -> '2.80'
$ws.Cells.Item(6,12).Value = "3.01"  # L6: '2.78' -> '3.01'

$ws.Cells.Item(7,4).Value = 52  # D7: 54 -> 52
$ws.Cells.Item(7,5).Value = "27,816,628.00"  # E7: '27,325,068.00' -> '27,816,628.00'
$ws.Cells.Item(7,6).Value = "360,218,661.00"  # F7: '347,149,837.00' -> '360,218,661.00'
$ws.Cells.Item(7,7).Value = "7.72"  # G7: '7.87' -> '7.72'
$ws.Cells.Item(7,8).Value = 10.685  # H7: 10.463 -> 10.685
$ws.Cells.Item(7,9).Value = 268  # I7: 254 -> 268
$ws.Cells.Item(7,11).Value = "4.09"  # K7: '4.34' -> '4.09'
$ws.Cells.Item(7,12).Value = "4.48"  # L7: '4.72' -> '4.48'

$ws.Cells.Item(8,4).Value = 43  # D8: 41 -> 43
$ws.Cells.Item(8,5).Value = "30,262,476.00"  # E8: '27,192,618.00' -> '30,262,476.00'
$ws.Cells.Item(8,6).Value = "342,367,394.00"  # F8: '334,624,168.00' -> '342,367,394.00'
$ws.Cells.Item(8,7).Value = "8.84"  # G8: '8.13' -> '8.84'
$ws.Cells.Item(8,8).Value = 21.142  # H8: 20.693 -> 21.142
$ws.Cells.Item(8,9).Value = 263  # I8: 249 -> 263
$ws.Cells.Item(8,10).Value = 13  # J8: 12 -> 13
$ws.Cells.Item(8,11).Value = "3.76"  # K8: '3.72' -> '3.76'
$ws.Cells.Item(8,12).Value = "4.94"  # L8: '4.82' -> '4.94'

$ws.Cells.Item(9,4).Value = 34  # D9: 38 -> 34
$ws.Cells.Item(9,5).Value = "30,316,641.00"  # E9: '29,079,025.00' -> '30,316,641.00'
$ws.Cells.Item(9,6).Value = "374,326,651.00"  # F9: '356,015,842.00' -> '374,326,651.00'
$ws.Cells.Item(9,7).Value = "8.10"  # G9: '8.17' -> '8.10'
$ws.Cells.Item(9,8).Value = 11.724  # H9: 11.396 -> 11.724
$ws.Cells.Item(9,9).Value = 270  # I9: 255 -> 270
$ws.Cells.Item(9,11).Value = "4.92"  # K9: '5.30' -> '4.92'
$ws.Cells.Item(9,12).Value = "2.59"  # L9: '2.75' -> '2.59'

$ws.Cells.Item(10,4).Value = 22  # D10: 21 -> 22
$ws.Cells.Item(10,5).Value = "15,618,789.00"  # E10: '15,497,789.00' -> '15,618,789.00'
$ws.Cells.Item(10,6).Value = "296,079,966.00"  # F10: '287,278,647.00' -> '296,079,966.00'
$ws.Cells.Item(10,7).Value = "5.28"  # G10: '5.39' -> '5.28'
$ws.Cells.Item(10,8).Value = 11.255  # H10: 11.12 -> 11.255
$ws.Cells.Item(10,9).Value = 203  # I10: 189 -> 203
$ws.Cells.Item(10,11).Value = "3.29"  # K10: '3.47' -> '3.29'
$ws.Cells.Item(10,12).Value = "2.96"  # L10: '3.17' -> '2.96'

$ws.Cells.Item(11,4).Value = 45  # D11: 51 -> 45
$ws.Cells.Item(11,6).Value = "341,063,822.00"  # F11: '327,671,563.00' -> '341,063,822.00'
$ws.Cells.Item(11,7).Value = "12.23"  # G11: '12.73' -> '12.23'
$ws.Cells.Item(11,8).Value = 12.179  # H11: 11.973 -> 12.179
$ws.Cells.Item(11,9).Value = 265  # I11: 250 -> 265
$ws.Cells.Item(11,11).Value = "9.76"  # K11: '10.37' -> '9.76'
$ws.Cells.Item(11,12).Value = "7.17"  # L11: '7.60' -> '7.17'

$ws.Cells.Item(12,4).Value = 45  # D12: 41 -> 45
$ws.Cells.Item(12,5).Value = "24,639,241.00"  # E12: '22,300,225.00' -> '24,639,241.00'
$ws.Cells.Item(12,6).Value = "332,107,238.00"  # F12: '315,752,427.00' -> '332,107,238.00'
$ws.Cells.Item(12,7).Value = "7.42"  # G12: '7.06' -> '7.42'
$ws.Cells.Item(12,8).Value = 17.115  # H12: 16.427 -> 17.115
$ws.Cells.Item(12,9).Value = 269  # I12: 255 -> 269
$ws.Cells.Item(12,11).Value = "5.67"  # K12: '6.11' -> '5.67'
$ws.Cells.Item(12,12).Value = "6.32"  # L12: '6.67' -> '6.32'

$ws.Cells.Item(13,4).Value = 34  # D13: 32 -> 34
$ws.Cells.Item(13,5).Value = "24,749,928.00"  # E13: '24,159,479.00' -> '24,749,928.00'
$ws.Cells.Item(13,6).Value = "348,362,716.00"  # F13: '339,780,179.00' -> '348,362,716.00'
$ws.Cells.Item(13,7).Value = "7.10"  # G13: '7.11' -> '7.10'
$ws.Cells.Item(13,8).Value = 15.169  # H13: 14.956 -> 15.169
$ws.Cells.Item(13,9).Value = 264  # I13: 250 -> 264
$ws.Cells.Item(13,11).Value = "3.47"  # K13: '3.62' -> '3.47'
$ws.Cells.Item(13,12).Value = "3.41"  # L13: '3.60' -> '3.41'

$ws.Cells.Item(14,4).Value = 39  # D14: 46 -> 39
$ws.Cells.Item(14,6).Value = "334,856,094.00"  # F14: '324,168,619.00' -> '334,856,094.00'
$ws.Cells.Item(14,7).Value = "12.03"  # G14: '12.43' -> '12.03'
$ws.Cells.Item(14,8).Value = 9.402  # H14: 9.268 -> 9.402
$ws.Cells.Item(14,9).Value = 264  # I14: 250 -> 264
$ws.Cells.Item(14,11).Value = "11.15"  # K14: '11.66' -> '11.15'
$ws.Cells.Item(14,12).Value = "4.92"  # L14: '5.20' -> '4.92'

$ws.Cells.Item(15,4).Value = 35  # D15: 32 -> 35
$ws.Cells.Item(15,5).Value = "25,731,270.00"  # E15: '24,333,873.00' -> '25,731,270.00'
$ws.Cells.Item(15,6).Value = "343,540,166.00"  # F15: '331,674,751.00' -> '343,540,166.00'
$ws.Cells.Item(15,7).Value = "7.49"  # G15: '7.34' -> '7.49'
$ws.Cells.Item(15,8).Value = 9.11  # H15: 8.732 -> 9.11
$ws.Cells.Item(15,9).Value = 268  # I15: 253 -> 268
$ws.Cells.Item(15,11).Value = "3.11"  # K15: '3.27' -> '3.11'
$ws.Cells.Item(15,12).Value = "3.73"  # L15: '3.95' -> '3.73'

$ws.Cells.Item(16,4).Value = 43  # D16: 42 -> 43
$ws.Cells.Item(16,5).Value = "30,304,175.00"  # E16: '29,648,492.00' -> '30,304,175.00'
$ws.Cells.Item(16,6).Value = "329,226,690.00"  # F16: '317,215,129.00' -> '329,226,690.00'
$ws.Cells.Item(16,7).Value = "9.20"  # G16: '9.35' -> '9.20'
$ws.Cells.Item(16,8).Value = 9.089  # H16: 8.795 -> 9.089
$ws.Cells.Item(16,9).Value = 268  # I16: 254 -> 268
$ws.Cells.Item(16,11).Value = "2.63"  # K16: '2.60' -> '2.63'
$ws.Cells.Item(16,12).Value = "4.10"  # L16: '4.33' -> '4.10'

$ws.Cells.Item(17,4).Value = 40  # D17: 37 -> 40
$ws.Cells.Item(17,5).Value = "30,628,998.00"  # E17: '29,957,114.00' -> '30,628,998.00'
$ws.Cells.Item(17,6).Value = "328,089,803.00"  # F17: '319,800,812.00' -> '328,089,803.00'
$ws.Cells.Item(17,7).Value = "9.34"  # G17: '9.37' -> '9.34'
$ws.Cells.Item(17,8).Value = 18.989  # H17: 18.604 -> 18.989
$ws.Cells.Item(17,9).Value = 266  # I17: 252 -> 266
$ws.Cells.Item(17,10).Value = 11  # J17: 9 -> 11
$ws.Cells.Item(17,11).Value = "2.98"  # K17: '2.93' -> '2.98'
$ws.Cells.Item(17,12).Value = "4.14"  # L17: '3.57' -> '4.14'

$ws.Cells.Item(18,4).Value = 27  # D18: 28 -> 27
$ws.Cells.Item(18,5).Value = "23,875,556.00"  # E18: '23,051,905.00' -> '23,875,556.00'
$ws.Cells.Item(18,6).Value = "281,143,423.00"  # F18: '271,676,597.00' -> '281,143,423.00'
$ws.Cells.Item(18,8).Value = 9.188  # H18: 8.915 -> 9.188
$ws.Cells.Item(18,9).Value = 169  # I18: 155 -> 169
$ws.Cells.Item(18,11).Value = "2.64"  # K18: '2.83' -> '2.64'
$ws.Cells.Item(18,12).Value = "1.78"  # L18: '1.94' -> '1.78'

